$d = $word.ActiveDocument

# Locate the unique "2018-05-08" submission-date text in the document.
$rng = $d.Content
$found = $rng.Find.Execute("2018-05-08", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # "2018-05-08" -> the digit to fix is the "0" right before the final "8"
    # ("2018-05-" is 8 characters long, so the "0" sits at offset 8).
    $zeroStart = $rng.Start + 8
    $zeroEnd = $zeroStart + 1

    $zeroRng = $d.Range($zeroStart, $zeroEnd)
    $zeroRng.Text = ""

    # Insert the replacement "1" as its own run (typed in place of the
    # deleted "0"), mirroring the edit 05-08 -> 05-18.
    $insPt = $d.Range($zeroStart, $zeroStart)
    $insPt.InsertBefore("1")

    # Force the freshly typed "1" onto its own run, distinct from the
    # surrounding "2018-05-" / "8" runs (toggling formatting on then back
    # off splits the run without leaving a visible formatting difference).
    $newRng = $d.Range($zeroStart, $zeroStart + 1)
    $newRng.Bold = $true
    $newRng.Bold = $false
}
